# "con loi tien hang" -- remove the stray reviewer comment on slide 4 and
# correct the cached "Update automatically" date shown on the slide
# master / every slide layout from 12/14/2019 to 12/18/2019.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Drop the lone comment (slide 4 / ppt/comments/comment1.xml).
# ---------------------------------------------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = $s.Comments.Count; $j -ge 1; $j--) {
        $s.Comments.Item($j).Delete()
    }
}

# ---------------------------------------------------------------------
# 2) Fix up the cached datetimeFigureOut text wherever it still reads
#    12/14/2019 -- the slide master and all seventeen slide layouts.
# ---------------------------------------------------------------------
$oldDate = "12/14/2019"
$newDate = "12/18/2019"
$ppPlaceholderDate = 16

function Update-DateShape($shape) {
    if (-not $shape.HasTextFrame) { return }
    if (-not $shape.TextFrame.HasText) { return }
    $isDatePlaceholder = $false
    if ($shape.Type -eq 14) {
        try {
            if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
    }
    if (-not $isDatePlaceholder) { return }
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq $oldDate) {
        $tr.Text = $newDate
    }
}

# Slide master.
$master = $p.SlideMaster
for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    Update-DateShape $master.Shapes.Item($si)
}

# Every slide layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        Update-DateShape $layout.Shapes.Item($si)
    }
}
